$d = $word.ActiveDocument

# --- Replace the "Predit" CA values (row 3, col 2) in the first four tables (2020-2023) ---
$d.Tables.Item(1).Cell(3, 2).Range.Text = "1165185"
$d.Tables.Item(2).Cell(3, 2).Range.Text = "1231697"
$d.Tables.Item(3).Cell(3, 2).Range.Text = "1299644"
$d.Tables.Item(4).Cell(3, 2).Range.Text = "1369025"

# --- Merge the split year-header runs ("202" + "N") into a single run ("202N") ---
# Tables 2-5 correspond to years 2021, 2022, 2023, 2024, whose header cell
# text is currently split across two separate runs ("202" and the last digit).
# Using Find/Replace on the cell range collapses the matched text into a
# single run carrying the full replacement string.
$years = @("2021", "2022", "2023", "2024")
for ($i = 0; $i -lt 4; $i++) {
    $tableIndex = $i + 2
    $headerCell = $d.Tables.Item($tableIndex).Cell(1, 1)
    $r = $headerCell.Range
    [void]$r.Find.Execute($years[$i], $true, $false, $false, $false, $false, $true, 1, $false, $years[$i], 2)
}
